# ---------------------------------------------------------------------------
# Adds a new worksheet "bras robotisés" (robotic arm component comparison)
# after "Feuil1", fills it in with the comparison table + a hyperlink, and
# lightly restyles/resizes the existing "Feuil1" sheet (vertical alignment of
# the data rows switches from "center" to "top", columns get wider, a few
# row heights change and the view is rezoomed/reselected).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Restyle "Feuil1" : the bordered data rows (2-5) move from vertical
#    "center" to vertical "top" alignment. Header row (row 1) is untouched.
# ---------------------------------------------------------------------------
$ws1.Range("A2:F5").VerticalAlignment = -4160   # xlTop

# Resize columns A:F to their new widths (best effort - the engine snaps to
# whole-pixel character widths, so we feed it the un-rounded target values).
$ws1.Columns.Item(1).ColumnWidth = 51.5546875
$ws1.Columns.Item(2).ColumnWidth = 54.6640625
$ws1.Columns.Item(3).ColumnWidth = 57.33203125
$ws1.Columns.Item(4).ColumnWidth = 57.44140625
$ws1.Columns.Item(5).ColumnWidth = 57.33203125
$ws1.Columns.Item(6).ColumnWidth = 53.33203125

# A few row heights changed slightly.
$ws1.Rows.Item(2).RowHeight = 150
$ws1.Rows.Item(3).RowHeight = 157.8
$ws1.Rows.Item(4).RowHeight = 130.2
$ws1.Rows.Item(5).RowHeight = 107.4

# ---------------------------------------------------------------------------
# 2) Insert the new worksheet right after "Feuil1".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "bras robotisés"

# ---------------------------------------------------------------------------
# 3) Fill in the comparison table. Columns A-C are filled in first (row by
#    row), and the "Lien" column D is added last - this mirrors the order in
#    which the original author entered the data (and keeps the shared
#    string table ordering identical).
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Critère"
$ws2.Range("B1").Value = "Valeur recommandée"
$ws2.Range("C1").Value = "Pourquoi c’est important"

$ws2.Range("A2").Value = "Nombre d’axes (degrés de liberté)"
$ws2.Range("B2").Value = "minimum 3–4 axes"
$ws2.Range("C2").Value = "suffisant pour poser une graine, plus pour des mouvements complexes"

$ws2.Range("A3").Value = "Portée / reach"
$ws2.Range("B3").Value = "~10–30 cm (ou plus selon taille serre)"
$ws2.Range("C3").Value = "pour atteindre différentes zones de plantation"

$ws2.Range("A4").Value = "Charge utile / payload"
$ws2.Range("B4").Value = "50–200 g (selon grain + mécanisme)"
$ws2.Range("C4").Value = "pour supporter le poids du mécanisme de semis"

$ws2.Range("A5").Value = "Précision / répétabilité"
$ws2.Range("B5").Value = "< ±1 cm (voire < ±5 mm)"
$ws2.Range("C5").Value = "pour positionner les graines correctement"

$ws2.Range("A6").Value = "Interface / contrôle"
$ws2.Range("B6").Value = "UART, PWM, I2C, librairie / API"
$ws2.Range("C6").Value = "pour l’intégrer avec ESP32 / Raspberry Pi"

$ws2.Range("A7").Value = "Matériau / robustesse"
$ws2.Range("B7").Value = "métal, aluminium, plastique renforcé"
$ws2.Range("C7").Value = "résister à l’humidité ou contraintes mécaniques"

$ws2.Range("A8").Value = "Support & documentation"
$ws2.Range("B8").Value = "code d’exemple, schémas, tutoriels"
$ws2.Range("C8").Value = "facilite l’intégration dans ton projet"

$ws2.Range("A9").Value = "Coût & disponibilité"
$ws2.Range("B9").Value = "budget faible / intermédiaire"
$ws2.Range("C9").Value = "pour ton prototype, pas besoin de bras industriel complet"

# Column D ("Lien") is added afterwards.
$ws2.Range("D1").Value = "Lien"
$ws2.Range("D2").Value = "https://www.lextronic.fr/bras-robotique-mearm-classic-maker-kit-63665.html"
$ws2.Range("D3").Value = "https://www.amazon.fr/robotique-m%C3%A9canique-Raspberry-%C3%A9tudiants-Bricolage/dp/B08P4WQ82H/ref=sr_1_24?crid=2NYO0WA18RQYS&dib=eyJ2IjoiMSJ9.5ENYusLGTOI7QWPhaq_mQmjmVKKSvUf24MaU_5X2Hv5ZxAWQ1ozNTjWIIwEys10Ba1AvODQwj-omkEoK8AFJb4M0_Y4wHypVprQgj2Vd5aNDUFE-kL_8mQXsDqFWqYf_aTLwWqCrisASQn2Zspo9vGxNxAE9YnsEo7ZJc5FXjk1-D1TI2zceO6-sGbKt7u1SZwbV707TG-uIMdRPgp9dk6KFcmeBthiooos8av9JVtHupYD_hd8mtc39MgVqRfUNc9NOF_aLElb04PEWYPLXXgDCW0uFyvGRXsGaye9-e_w.TVI2qY_OnTWRlt5ADNceLDMtlvh9nPRF-AJGI-tjUAg&dib_tag=se&keywords=bras+robotique&qid=1759327465&sprefix=bras+r%2Caps%2C875&sr=8-24"
$ws2.Range("D4").Value = "https://www.amazon.fr/Tatiy-Programmation-Servomoteurs-lenseignement-Universitaire/dp/B0B4VKYQ6S/ref=sr_1_29?crid=2NYO0WA18RQYS&dib=eyJ2IjoiMSJ9.5ENYusLGTOI7QWPhaq_mQmjmVKKSvUf24MaU_5X2Hv5ZxAWQ1ozNTjWIIwEys10Ba1AvODQwj-omkEoK8AFJb4M0_Y4wHypVprQgj2Vd5aNDUFE-kL_8mQXsDqFWqYf_aTLwWqCrisASQn2Zspo9vGxNxAE9YnsEo7ZJc5FXjk1-D1TI2zceO6-sGbKt7u1SZwbV707TG-uIMdRPgp9dk6KFcmeBthiooos8av9JVtHupYD_hd8mtc39MgVqRfUNc9NOF_aLElb04PEWYPLXXgDCW0uFyvGRXsGaye9-e_w.TVI2qY_OnTWRlt5ADNceLDMtlvh9nPRF-AJGI-tjUAg&dib_tag=se&keywords=bras+robotique&qid=1759327465&sprefix=bras+r%2Caps%2C875&sr=8-29"

# Turn D2 into a real hyperlink (this also applies the built-in "Hyperlink"
# style: underlined, themed colour).
$ws2.Hyperlinks.Add($ws2.Range("D2"), $ws2.Range("D2").Value) | Out-Null

# ---------------------------------------------------------------------------
# 4) Formatting.
# ---------------------------------------------------------------------------
# Header row: bold, centered, wrapped.
$header = $ws2.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4108     # xlCenter
$header.WrapText = $true

# Column A (criteria names): bold, vertically centered, wrapped.
$colA = $ws2.Range("A2:A9")
$colA.Font.Bold = $true
$colA.VerticalAlignment = -4108       # xlCenter
$colA.WrapText = $true

# Columns B and C: regular weight, vertically centered, wrapped.
$colBC = $ws2.Range("B2:C9")
$colBC.VerticalAlignment = -4108      # xlCenter
$colBC.WrapText = $true

# ---------------------------------------------------------------------------
# 5) Column widths / row heights.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 26.44140625
$ws2.Columns.Item(2).ColumnWidth = 41.77734375
$ws2.Columns.Item(3).ColumnWidth = 44.109375
$ws2.Columns.Item(4).ColumnWidth = 111.6640625

$ws2.Rows.Item(2).RowHeight = 43.2
$ws2.Rows.Item(3).RowHeight = 28.8
$ws2.Rows.Item(4).RowHeight = 28.8
$ws2.Rows.Item(5).RowHeight = 28.8
$ws2.Rows.Item(6).RowHeight = 28.8
$ws2.Rows.Item(7).RowHeight = 28.8
$ws2.Rows.Item(8).RowHeight = 28.8
$ws2.Rows.Item(9).RowHeight = 43.2
$ws2.Rows.Item(10).RowHeight = 33.6

# ---------------------------------------------------------------------------
# 6) View state: Feuil1 zoomed out to 58% with D4 selected (no longer the
#    active tab); the new sheet is active, zoomed at 100%, with D7 selected.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 58
$ws1.Range("D4").Select()

$ws2.Activate()
$ws2.Range("D7").Select()
